$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 220
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = ""
$ws.Range("H138").Value = 2303.0754
$ws.Range("J138").Value = 2605.738
$ws.Range("L138").Value = 7817.214
$ws.Range("N138").Value = -18097.214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 8366.666999999999
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 12
$ws.Range("H32").Value = 138625.14
$ws.Range("I32").Value = 149484.02
$ws.Range("K32").Value = 149484.02
$ws.Range("M32").Value = -149197.02
$ws.Range("H74").Value = 3986860.2
$ws.Range("I74").Value = 13891139
$ws.Range("J74").Value = 25148.8
$ws.Range("K74").Value = 13891139
$ws.Range("L74").Value = 25148.8
$ws.Range("M74").Value = -13890265
$ws.Range("N74").Value = -26896.8
$ws.Range("H77").Value = 3986860.2
$ws.Range("I77").Value = 13891139
$ws.Range("J77").Value = 25148.8
$ws.Range("K77").Value = 69455695
$ws.Range("L77").Value = 125744
$ws.Range("M77").Value = -69451327
$ws.Range("N77").Value = -134480
$ws.Range("H110").Value = 1256.4117
$ws.Range("I110").Value = 1028
$ws.Range("K110").Value = 1028
$ws.Range("M110").Value = 1017

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 8366.666999999999
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 15
$ws.Range("H94").Value = 5979.1113
$ws.Range("I94").Value = 5105.923
$ws.Range("K94").Value = 5105.923
$ws.Range("M94").Value = -4654.923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1427112.8
$ws.Range("J31").Value = 2869.9285
$ws.Range("L31").Value = 2869.9285
$ws.Range("N31").Value = -3459.9285
$ws.Range("H34").Value = 1427112.8
$ws.Range("J34").Value = 2869.9285
$ws.Range("L34").Value = 2869.9285
$ws.Range("N34").Value = -3273.9285
$ws.Range("H58").Value = 5963488.5
$ws.Range("I58").Value = 3676.75
$ws.Range("K58").Value = 3676.75
$ws.Range("M58").Value = -3473.75
$ws.Range("H99").Value = 18100
$ws.Range("I99").Value = 20450.25
$ws.Range("K99").Value = 20450.25
$ws.Range("M99").Value = -18952.25
$ws.Range("H126").Value = 18100
$ws.Range("I126").Value = 20450.25
$ws.Range("K126").Value = 61350.75
$ws.Range("M126").Value = -58880.75
$ws.Range("H136").Value = 5963488.5
$ws.Range("I136").Value = 3676.75
$ws.Range("K136").Value = 11030.25
$ws.Range("M136").Value = -8480.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 348.0909
$ws.Range("J92").Value = 313
$ws.Range("L92").Value = 939
$ws.Range("N92").Value = -3435
$ws.Range("H127").Value = 10833.333
$ws.Range("J127").Value = 10833.333
$ws.Range("L127").Value = 32499.999
$ws.Range("N127").Value = -42419.999
$ws.Range("H129").Value = 1392.8334
$ws.Range("I129").Value = 737.6667
$ws.Range("J129").Value = 3358.3333
$ws.Range("K129").Value = 2213.0001
$ws.Range("L129").Value = 10074.9999
$ws.Range("M129").Value = 2786.9999
$ws.Range("N129").Value = -20074.9999
$ws.Range("H131").Value = 7560.1113
$ws.Range("J131").Value = 9337.429
$ws.Range("L131").Value = 28012.287
$ws.Range("N131").Value = -38092.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 28593.334
$ws.Range("J63").Value = 28593.334
$ws.Range("L63").Value = 28593.334
$ws.Range("N63").Value = -29965.334
$ws.Range("H66").Value = 28593.334
$ws.Range("J66").Value = 28593.334
$ws.Range("L66").Value = 85780.00199999999
$ws.Range("N66").Value = -92644.00199999999
$ws.Range("H98").Value = 10045
$ws.Range("J98").Value = 10045
$ws.Range("L98").Value = 10045
$ws.Range("N98").Value = -16035
$ws.Range("H102").Value = 2983.2856
$ws.Range("I102").Value = 2857.5
$ws.Range("K102").Value = 2857.5
$ws.Range("M102").Value = -1235.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4138.647
$ws.Range("I7").Value = 3964
$ws.Range("J7").Value = 4706.25
$ws.Range("K7").Value = 3964
$ws.Range("L7").Value = 4706.25
$ws.Range("M7").Value = -3852
$ws.Range("N7").Value = -4930.25
$ws.Range("H26").Value = 15000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 15000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = -15590
$ws.Range("H31").Value = 3116.5
$ws.Range("I31").Value = 329
$ws.Range("J31").Value = 7762.3335
$ws.Range("K31").Value = 329
$ws.Range("L31").Value = 7762.3335
$ws.Range("M31").Value = -81
$ws.Range("N31").Value = -8258.333500000001
$ws.Range("H46").Value = 4012.125
$ws.Range("J46").Value = 4513.857
$ws.Range("L46").Value = 4513.857
$ws.Range("N46").Value = -4889.857
$ws.Range("H100").Value = 5195.1113
$ws.Range("I100").Value = 5074.2
$ws.Range("K100").Value = 5074.2
$ws.Range("M100").Value = -4533.2
$ws.Range("H122").Value = 2590.2727
$ws.Range("I122").Value = 2349.3
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7047.900000000001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4597.900000000001
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 4138.647
$ws.Range("I126").Value = 3964
$ws.Range("J126").Value = 4706.25
$ws.Range("K126").Value = 11892
$ws.Range("L126").Value = 14118.75
$ws.Range("M126").Value = -9422
$ws.Range("N126").Value = -19058.75
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 126900.78
$ws.Range("I122").Value = 1998.2
$ws.Range("K122").Value = 5994.6
$ws.Range("M122").Value = -3544.6
$ws.Range("H136").Value = 121739130
$ws.Range("I136").Value = 43478260
$ws.Range("J136").Value = 200000000
$ws.Range("K136").Value = 130434780
$ws.Range("L136").Value = 600000000
$ws.Range("M136").Value = -130432230
$ws.Range("N136").Value = -600005100
